$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") values
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1413
$ws.Range("F4").Value = 19829
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 1093
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 504
$ws.Range("F11").Value = 729
$ws.Range("F12").Value = 258
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 110
$ws.Range("F17").Value = 232
$ws.Range("F18").Value = 191
$ws.Range("F19").Value = 1333
$ws.Range("F20").Value = 0
$ws.Range("F23").Value = 46
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 66
$ws.Range("F26").Value = 318
$ws.Range("F27").Value = 1092
$ws.Range("F33").Value = 57
$ws.Range("F34").Value = 0
$ws.Range("F38").Value = 12574
$ws.Range("F39").Value = 1329
$ws.Range("F41").Value = 0
$ws.Range("F43").Value = 255
$ws.Range("F44").Value = 353
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 93

# Sheet "演出" (sheet2): update column F value
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0

# Sheet "全部类型" (sheet4): update column F values
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 19829
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 1093
$ws.Range("F9").Value = 7493
$ws.Range("F10").Value = 504
$ws.Range("F12").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 1333
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 1092
$ws.Range("F28").Value = 26
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 5220
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 57
$ws.Range("F36").Value = 2806
$ws.Range("F39").Value = 50
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 1329
$ws.Range("F44").Value = 54
$ws.Range("F45").Value = 255
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 318
$ws.Range("F49").Value = 0
